# Update construct_texts worksheet:
#  - Replace the "percent free lunch qualified" construct with the
#    2019-20 NCES "percent free and reduced lunch" construct
#  - Add new rows for the "number of strategies" (ss1..ss5) constructs
#    and the school-level dummy variables (schoollevelHS / schoollevelMS)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows: number-of-strategies constructs (column A, then column B) ---
$ws.Range("A27").Value = "ss1"
$ws.Range("A28").Value = "ss2"
$ws.Range("A29").Value = "ss3"
$ws.Range("A30").Value = "ss4"
$ws.Range("A31").Value = "ss5"

$ws.Range("B27").Value = "One strategy"
$ws.Range("B28").Value = "Two strategies"
$ws.Range("B29").Value = "Three strategies"
$ws.Range("B30").Value = "Four strategies"
$ws.Range("B31").Value = "Five strategies"

# --- New rows: school level dummy variables (column A, then column B) ---
$ws.Range("A32").Value = "schoollevelHS"
$ws.Range("A33").Value = "schoollevelMS"

$ws.Range("B32").Value = "High school"
$ws.Range("B33").Value = "Middle school"

# --- Replace free lunch construct with free & reduced lunch construct ---
$ws.Range("B21").Value = "Percent free and reduced lunch"
$ws.Range("A21").Value = "percentfreereducedlunch"

# --- Update selection shown in the sheet view ---
$ws.Activate()
$ws.Range("A21").Select()
